$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing header/value cells
$ws.Range("C1").Value = "data_sorteio"
$ws.Range("B2").Value = "Consulta Odontológica Bonificada"
$ws.Range("C2").Value = "08/08/2025 16:36:11"

# Add new column D: header (copy header formatting) + value
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D1").Value = "atendente_cadastro"

$ws.Range("D2").Value = "Lucas Mendes"
